$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 812, shifting existing rows 812:853 down to 813:854
$ws.Rows.Item(812).Insert()

# Populate the newly inserted row 812 with the new data point.
# Column A holds a date-formatted string (e.g. "2026/02/15"); without the
# leading apostrophe Excel auto-converts it to a date serial number, which
# doesn't match the source data (plain text dates, like the rest of the
# column). The apostrophe forces text entry, then resetting the style back
# to "Normal" clears the "Text" number-format style COM applies when it
# detects the forced-text entry, so the cell ends up with no explicit style
# -- same as its sibling cells.
$ws.Cells.Item(812, 1).Value = "'2026/02/15"
$ws.Cells.Item(812, 1).Style = "Normal"
$ws.Cells.Item(812, 2).Value = "日"
$ws.Cells.Item(812, 3).Value = 4
$ws.Cells.Item(812, 4).Value = 201
